# Auto-generated edit script applying the Ravana_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2764.5715
$ws.Range("I19").Value = 2900.3333
$ws.Range("J19").Value = 1950
$ws.Range("K19").Value = 2900.3333
$ws.Range("L19").Value = 1950
$ws.Range("M19").Value = -2725.3333
$ws.Range("N19").Value = -2300

$ws.Range("H45").Value = 8766.666999999999
$ws.Range("I45").Value = 2875
$ws.Range("J45").Value = 20550
$ws.Range("K45").Value = 8625
$ws.Range("L45").Value = 61650
$ws.Range("M45").Value = -8433
$ws.Range("N45").Value = -62034

$ws.Range("H58").Value = 7480
$ws.Range("I58").Value = 450
$ws.Range("J58").Value = 12166.667
$ws.Range("K58").Value = 1350
$ws.Range("L58").Value = 36500.001
$ws.Range("M58").Value = -1200
$ws.Range("N58").Value = -36800.001

$ws.Range("H80").Value = 705.2727
$ws.Range("I80").Value = 702.5714
$ws.Range("J80").Value = 710
$ws.Range("K80").Value = 2107.7142
$ws.Range("L80").Value = 2130
$ws.Range("M80").Value = -1109.7142
$ws.Range("N80").Value = -4126

$ws.Range("H83").Value = 705.2727
$ws.Range("I83").Value = 702.5714
$ws.Range("J83").Value = 710
$ws.Range("K83").Value = 6323.1426
$ws.Range("L83").Value = 6390
$ws.Range("M83").Value = -1331.1426
$ws.Range("N83").Value = -16374

$ws.Range("H100").Value = 7798.2
$ws.Range("I100").Value = 2994
$ws.Range("K100").Value = 2994
$ws.Range("M100").Value = -2453

$ws.Range("H112").Value = 1772.1875
$ws.Range("J112").Value = 1958.0769
$ws.Range("L112").Value = 5874.2307
$ws.Range("N112").Value = -8090.2307

$ws.Range("H113").Value = 3495
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H132").Value = 1177.9286
$ws.Range("I132").Value = 1177.9286
$ws.Range("K132").Value = 3533.7858
$ws.Range("M132").Value = -1003.7858

$ws.Range("H137").Value = 2477.3667
$ws.Range("I137").Value = 1366
$ws.Range("K137").Value = 4098
$ws.Range("M137").Value = -1548

$ws.Range("H138").Value = 4140.516
$ws.Range("J138").Value = 4272.593
$ws.Range("L138").Value = 12817.779
$ws.Range("N138").Value = -23097.779

$ws.Range("H141").Value = 8496.75
$ws.Range("I141").Value = 7995.6665
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 23986.9995
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -18806.9995
$ws.Range("N141").Value = -40360


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 3000
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = -2670

$ws.Range("H32").Value = 6202.4365
$ws.Range("I32").Value = 5946.926
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 5946.926
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -5659.926
$ws.Range("N32").Value = -20574

$ws.Range("H45").Value = 2611.75
$ws.Range("I45").Value = 2611.75
$ws.Range("K45").Value = 2611.75
$ws.Range("M45").Value = -2234.75

$ws.Range("H74").Value = 1907.2727
$ws.Range("I74").Value = 1531.1111
$ws.Range("K74").Value = 1531.1111
$ws.Range("M74").Value = -657.1111000000001

$ws.Range("H77").Value = 1907.2727
$ws.Range("I77").Value = 1531.1111
$ws.Range("K77").Value = 7655.5555
$ws.Range("M77").Value = -3287.5555

$ws.Range("H110").Value = 1133
$ws.Range("I110").Value = 1000
$ws.Range("J110").Value = 1199.5
$ws.Range("K110").Value = 1000
$ws.Range("L110").Value = 1199.5
$ws.Range("M110").Value = 1045
$ws.Range("N110").Value = -5289.5

$ws.Range("H122").Value = 3063
$ws.Range("I122").Value = 3032.7778
$ws.Range("J122").Value = 3199
$ws.Range("K122").Value = 9098.3334
$ws.Range("L122").Value = 9597
$ws.Range("M122").Value = -6648.3334
$ws.Range("N122").Value = -14497


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 677.6
$ws.Range("I22").Value = 677
$ws.Range("K22").Value = 677
$ws.Range("M22").Value = -327

$ws.Range("H31").Value = 3303.5
$ws.Range("J31").Value = 5064.4287
$ws.Range("L31").Value = 5064.4287
$ws.Range("N31").Value = -5654.4287

$ws.Range("H34").Value = 3303.5
$ws.Range("J34").Value = 5064.4287
$ws.Range("L34").Value = 5064.4287
$ws.Range("N34").Value = -5468.4287

$ws.Range("H99").Value = 10000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H107").Value = 2474.182
$ws.Range("J107").Value = 2816.2222
$ws.Range("L107").Value = 2816.2222
$ws.Range("N107").Value = -6656.2222

$ws.Range("H122").Value = 2359.6
$ws.Range("I122").Value = 599.3333
$ws.Range("K122").Value = 1797.9999
$ws.Range("M122").Value = 652.0001

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H134").Value = 3833
$ws.Range("J134").Value = 3499
$ws.Range("L134").Value = 10497
$ws.Range("N134").Value = -15567


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25815398
$ws.Range("I4").Value = 25126946
$ws.Range("K4").Value = 75380838
$ws.Range("M4").Value = -75380726

$ws.Range("H92").Value = 1173.75
$ws.Range("J92").Value = 1500
$ws.Range("L92").Value = 4500
$ws.Range("N92").Value = -6996

$ws.Range("H107").Value = 282.25
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 282.25
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 846.75
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4686.75

$ws.Range("H113").Value = 2641.7144
$ws.Range("J113").Value = 2641.7144
$ws.Range("L113").Value = 7925.1432
$ws.Range("N113").Value = -12265.1432

$ws.Range("H140").Value = 1481.875
$ws.Range("I140").Value = 1481.875
$ws.Range("K140").Value = 4445.625
$ws.Range("M140").Value = 734.375


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5471.857
$ws.Range("I80").Value = 4716.1665
$ws.Range("K80").Value = 4716.1665
$ws.Range("M80").Value = -3718.1665

$ws.Range("H83").Value = 5471.857
$ws.Range("I83").Value = 4716.1665
$ws.Range("K83").Value = 23580.8325
$ws.Range("M83").Value = -18588.8325


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5027
$ws.Range("I22").Value = 3706
$ws.Range("J22").Value = 5687.5
$ws.Range("K22").Value = 3706
$ws.Range("L22").Value = 5687.5
$ws.Range("M22").Value = -3411
$ws.Range("N22").Value = -6277.5

$ws.Range("H27").Value = 5027
$ws.Range("I27").Value = 3706
$ws.Range("J27").Value = 5687.5
$ws.Range("K27").Value = 3706
$ws.Range("L27").Value = 5687.5
$ws.Range("M27").Value = -3599
$ws.Range("N27").Value = -5901.5

$ws.Range("H82").Value = 2903.5
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 2903.5
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H132").Value = 5186.25
$ws.Range("I132").Value = 4915.3335
$ws.Range("K132").Value = 14746.0005
$ws.Range("M132").Value = -12216.0005


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 19833.334
$ws.Range("J54").Value = 19833.334
$ws.Range("L54").Value = 19833.334
$ws.Range("N54").Value = -20873.334

$ws.Range("H126").Value = 1905.5
$ws.Range("I126").Value = 1905.5
$ws.Range("K126").Value = 5716.5
$ws.Range("M126").Value = -3246.5

$ws.Range("H132").Value = 3180.7083
$ws.Range("I132").Value = 2375.3572
$ws.Range("J132").Value = 4308.2
$ws.Range("K132").Value = 7126.071599999999
$ws.Range("L132").Value = 12924.6
$ws.Range("M132").Value = -4596.071599999999
$ws.Range("N132").Value = -17984.6

$ws.Range("H136").Value = 3004.4075
$ws.Range("I136").Value = 2816.95
$ws.Range("J136").Value = 3540
$ws.Range("K136").Value = 8450.849999999999
$ws.Range("L136").Value = 10620
$ws.Range("M136").Value = -5900.849999999999
$ws.Range("N136").Value = -15720


Write-Output "edits applied"